$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1061.7094
$ws.Cells.Item(17, 10).Value = 1130.6072
$ws.Cells.Item(17, 12).Value = 3391.8216
$ws.Cells.Item(17, 14).Value = -3727.8216
$ws.Cells.Item(33, 8).Value = 1020.2222
$ws.Cells.Item(33, 9).Value = 1085.25
$ws.Cells.Item(33, 11).Value = 1085.25
$ws.Cells.Item(33, 13).Value = -856.25
$ws.Cells.Item(43, 8).Value = 1642.7037
$ws.Cells.Item(43, 9).Value = 1492.5
$ws.Cells.Item(43, 10).Value = 1705.9474
$ws.Cells.Item(43, 11).Value = 1492.5
$ws.Cells.Item(43, 12).Value = 1705.9474
$ws.Cells.Item(43, 13).Value = -1423.5
$ws.Cells.Item(43, 14).Value = -1843.9474
$ws.Cells.Item(74, 8).Value = 3495.3
$ws.Cells.Item(77, 8).Value = 3495.3
$ws.Cells.Item(112, 8).Value = 5060.095
$ws.Cells.Item(112, 10).Value = 5264.1
$ws.Cells.Item(112, 12).Value = 15792.3
$ws.Cells.Item(112, 14).Value = -18008.3
$ws.Cells.Item(140, 8).Value = 74398.57000000001
$ws.Cells.Item(140, 10).Value = 74398.57000000001
$ws.Cells.Item(140, 12).Value = 74398.57000000001
$ws.Cells.Item(140, 14).Value = -84758.57000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 34091.168
$ws.Cells.Item(7, 9).Value = 20323.5
$ws.Cells.Item(7, 10).Value = 40975
$ws.Cells.Item(7, 11).Value = 20323.5
$ws.Cells.Item(7, 12).Value = 40975
$ws.Cells.Item(7, 13).Value = -20209.5
$ws.Cells.Item(7, 14).Value = -41203
$ws.Cells.Item(32, 8).Value = 8782632
$ws.Cells.Item(32, 9).Value = 9269167
$ws.Cells.Item(32, 10).Value = 25000
$ws.Cells.Item(32, 11).Value = 9269167
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = -9268880
$ws.Cells.Item(32, 14).Value = -25574
$ws.Cells.Item(45, 8).Value = 2511.647
$ws.Cells.Item(45, 9).Value = 2051.3
$ws.Cells.Item(45, 11).Value = 2051.3
$ws.Cells.Item(45, 13).Value = -1674.3
$ws.Cells.Item(52, 8).Value = 92520
$ws.Cells.Item(52, 10).Value = 92520
$ws.Cells.Item(52, 12).Value = 92520
$ws.Cells.Item(52, 14).Value = -93156
$ws.Cells.Item(101, 8).Value = 76447.17999999999
$ws.Cells.Item(101, 10).Value = 76447.17999999999
$ws.Cells.Item(101, 12).Value = 76447.17999999999
$ws.Cells.Item(101, 14).Value = -82937.17999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 88712
$ws.Cells.Item(2, 10).Value = 88712
$ws.Cells.Item(2, 12).Value = 88712
$ws.Cells.Item(2, 14).Value = -88938
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).ClearContents()
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 20107
$ws.Cells.Item(61, 10).Value = 20107
$ws.Cells.Item(61, 12).Value = 20107
$ws.Cells.Item(61, 14).Value = -20733
$ws.Cells.Item(98, 8).Value = 79771
$ws.Cells.Item(98, 10).Value = 79771
$ws.Cells.Item(98, 12).Value = 79771
$ws.Cells.Item(98, 14).Value = -85761
$ws.Cells.Item(100, 8).Value = 84821.5
$ws.Cells.Item(100, 10).Value = 84821.5
$ws.Cells.Item(100, 12).Value = 84821.5
$ws.Cells.Item(100, 14).Value = -86985.5
$ws.Cells.Item(118, 8).Value = 47952
$ws.Cells.Item(118, 10).Value = 47952
$ws.Cells.Item(118, 12).Value = 47952
$ws.Cells.Item(118, 14).Value = -51266
$ws.Cells.Item(125, 8).Value = 98780
$ws.Cells.Item(125, 10).Value = 98780
$ws.Cells.Item(125, 12).Value = 98780
$ws.Cells.Item(125, 13).Value = -108620

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 25225
$ws.Cells.Item(12, 9).Value = 300
$ws.Cells.Item(12, 10).Value = 100000
$ws.Cells.Item(12, 11).Value = 300
$ws.Cells.Item(12, 12).Value = 100000
$ws.Cells.Item(12, 13).Value = -130
$ws.Cells.Item(12, 14).Value = -100340
$ws.Cells.Item(31, 8).Value = 5104.4404
$ws.Cells.Item(31, 9).Value = 1636.36
$ws.Cells.Item(31, 10).Value = 6573.9663
$ws.Cells.Item(31, 11).Value = 1636.36
$ws.Cells.Item(31, 12).Value = 6573.9663
$ws.Cells.Item(31, 13).Value = -1341.36
$ws.Cells.Item(31, 14).Value = -7163.9663
$ws.Cells.Item(34, 8).Value = 5104.4404
$ws.Cells.Item(34, 9).Value = 1636.36
$ws.Cells.Item(34, 10).Value = 6573.9663
$ws.Cells.Item(34, 11).Value = 1636.36
$ws.Cells.Item(34, 12).Value = 6573.9663
$ws.Cells.Item(34, 13).Value = -1434.36
$ws.Cells.Item(34, 14).Value = -6977.9663
$ws.Cells.Item(106, 8).Value = 43723.668
$ws.Cells.Item(106, 10).Value = 43723.668
$ws.Cells.Item(106, 12).Value = 43723.668
$ws.Cells.Item(106, 14).Value = -46247.668
$ws.Cells.Item(107, 8).Value = 5682927.5
$ws.Cells.Item(107, 9).Value = 10417482
$ws.Cells.Item(107, 10).Value = 1462.6
$ws.Cells.Item(107, 11).Value = 10417482
$ws.Cells.Item(107, 12).Value = 1462.6
$ws.Cells.Item(107, 13).Value = -10415562
$ws.Cells.Item(107, 14).Value = -5302.6

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 144.66667
$ws.Cells.Item(2, 9).Value = 39.666668
$ws.Cells.Item(2, 10).Value = 170.91667
$ws.Cells.Item(2, 11).Value = 238.000008
$ws.Cells.Item(2, 12).Value = 1025.50002
$ws.Cells.Item(2, 13).Value = -125.000008
$ws.Cells.Item(2, 14).Value = -1251.50002
$ws.Cells.Item(113, 8).Value = 641.9474
$ws.Cells.Item(113, 9).Value = 868
$ws.Cells.Item(113, 10).Value = 561.2143
$ws.Cells.Item(113, 11).Value = 2604
$ws.Cells.Item(113, 12).Value = 1683.6429
$ws.Cells.Item(113, 13).Value = -434
$ws.Cells.Item(113, 14).Value = -6023.6429
$ws.Cells.Item(139, 8).Value = 289599.38
$ws.Cells.Item(139, 9).Value = 528330
$ws.Cells.Item(139, 10).Value = 6106.75
$ws.Cells.Item(139, 11).Value = 1584990
$ws.Cells.Item(139, 12).Value = 18320.25
$ws.Cells.Item(139, 13).Value = -1579850
$ws.Cells.Item(139, 14).Value = -28600.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 10456
$ws.Cells.Item(99, 9).Value = 3820.3076
$ws.Cells.Item(99, 11).Value = 3820.3076
$ws.Cells.Item(99, 13).Value = -1574.3076
$ws.Cells.Item(132, 8).Value = 26321066
$ws.Cells.Item(132, 9).Value = 41672908
$ws.Cells.Item(132, 10).Value = 3618.1428
$ws.Cells.Item(132, 11).Value = 125018724
$ws.Cells.Item(132, 12).Value = 10854.4284
$ws.Cells.Item(132, 13).Value = -125016194
$ws.Cells.Item(132, 14).Value = -15914.4284

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4870.731
$ws.Cells.Item(7, 9).Value = 4675.643
$ws.Cells.Item(7, 10).Value = 5098.3335
$ws.Cells.Item(7, 11).Value = 4675.643
$ws.Cells.Item(7, 12).Value = 5098.3335
$ws.Cells.Item(7, 13).Value = -4563.643
$ws.Cells.Item(7, 14).Value = -5322.3335
$ws.Cells.Item(14, 8).Value = 82502.5
$ws.Cells.Item(14, 10).Value = 82502.5
$ws.Cells.Item(14, 12).Value = 82502.5
$ws.Cells.Item(14, 14).Value = -82846.5
$ws.Cells.Item(22, 8).Value = 20820.6
$ws.Cells.Item(22, 9).Value = 1050
$ws.Cells.Item(22, 10).Value = 34001
$ws.Cells.Item(22, 11).Value = 1050
$ws.Cells.Item(22, 12).Value = 34001
$ws.Cells.Item(22, 13).Value = -755
$ws.Cells.Item(22, 14).Value = -34591
$ws.Cells.Item(27, 8).Value = 20820.6
$ws.Cells.Item(27, 9).Value = 1050
$ws.Cells.Item(27, 10).Value = 34001
$ws.Cells.Item(27, 11).Value = 1050
$ws.Cells.Item(27, 12).Value = 34001
$ws.Cells.Item(27, 13).Value = -943
$ws.Cells.Item(27, 14).Value = -34215
$ws.Cells.Item(100, 8).Value = 78462.25
$ws.Cells.Item(100, 9).Value = 115652.875
$ws.Cells.Item(100, 10).Value = 4081
$ws.Cells.Item(100, 11).Value = 115652.875
$ws.Cells.Item(100, 12).Value = 4081
$ws.Cells.Item(100, 13).Value = -115111.875
$ws.Cells.Item(100, 14).Value = -5163
$ws.Cells.Item(125, 8).Value = 98476.664
$ws.Cells.Item(125, 10).Value = 98476.664
$ws.Cells.Item(125, 12).Value = 98476.664
$ws.Cells.Item(125, 14).Value = -108316.664
$ws.Cells.Item(126, 8).Value = 4870.731
$ws.Cells.Item(126, 9).Value = 4675.643
$ws.Cells.Item(126, 10).Value = 5098.3335
$ws.Cells.Item(126, 11).Value = 14026.929
$ws.Cells.Item(126, 12).Value = 15295.0005
$ws.Cells.Item(126, 13).Value = -11556.929
$ws.Cells.Item(126, 14).Value = -20235.0005
$ws.Cells.Item(132, 8).Value = 3303.276
$ws.Cells.Item(132, 9).Value = 2933.1667
$ws.Cells.Item(132, 10).Value = 3908.9092
$ws.Cells.Item(132, 11).Value = 8799.500100000001
$ws.Cells.Item(132, 12).Value = 11726.7276
$ws.Cells.Item(132, 13).Value = -6269.500100000001
$ws.Cells.Item(132, 14).Value = -16786.7276

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3437.111
$ws.Cells.Item(81, 9).Value = 3772.5
$ws.Cells.Item(81, 10).Value = 2766.3333
$ws.Cells.Item(81, 11).Value = 7545
$ws.Cells.Item(81, 12).Value = 5532.6666
$ws.Cells.Item(81, 13).Value = -6484
$ws.Cells.Item(81, 14).Value = -7654.6666
$ws.Cells.Item(84, 8).Value = 3437.111
$ws.Cells.Item(84, 9).Value = 3772.5
$ws.Cells.Item(84, 10).Value = 2766.3333
$ws.Cells.Item(84, 11).Value = 37725
$ws.Cells.Item(84, 12).Value = 27663.333
$ws.Cells.Item(84, 13).Value = -32421
$ws.Cells.Item(84, 14).Value = -38271.333
$ws.Cells.Item(109, 8).Value = 4020188.5
$ws.Cells.Item(109, 10).Value = 40377
$ws.Cells.Item(109, 12).Value = 40377
$ws.Cells.Item(109, 14).Value = -43151
$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 13).ClearContents()
$ws.Cells.Item(123, 14).ClearContents()
